# Updates cryptos list values (Price / Volume(1h)) per scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.374.71"
$ws.Range("E2").Value = "  -2.85%  "
$ws.Range("D3").Value = "1.740.94"
$ws.Range("E3").Value = "  -3.60%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "322.74"
$ws.Range("E5").Value = "  -3.89%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.4226"
$ws.Range("E7").Value = "  -9.92%  "
$ws.Range("D8").Value = "0.3609"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.50"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").Value = "0.07418"
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("D11").Value = "1.112"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  -5.06%  "
$ws.Range("D14").Value = "6.071"
$ws.Range("E14").Value = "  -4.96%  "
$ws.Range("D15").Value = "7.183"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").Value = "1.735.99"
$ws.Range("E16").Value = "  -3.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001063"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("D18").Value = "87.61"
$ws.Range("E18").Value = "  +5.83%  "
$ws.Range("D19").Value = "0.05996"
$ws.Range("E19").Value = "  -10.98%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "16.81"
$ws.Range("E21").Value = "  -4.00%  "
$ws.Range("D22").Value = "6.097"
$ws.Range("E22").Value = "  -5.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5230"
$ws.Range("E23").Value = "  -5.95%  "
$ws.Range("D24").Value = "27.402.60"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "11.36"
$ws.Range("E25").Value = "  -4.95%  "
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").Value = "20.15"
$ws.Range("E27").Value = "  -4.07%  "
$ws.Range("D28").Value = "2.376"
$ws.Range("E28").Value = "  -1.44%  "
$ws.Range("D29").Value = "149.24"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("D30").Value = "1.930.61"
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("D31").Value = "126.57"
$ws.Range("E31").Value = "  -5.99%  "
$ws.Range("E32").Value = "  -8.09%  "
$ws.Range("D33").Value = "5.686"
$ws.Range("D34").Value = "0.09105"
$ws.Range("E34").Value = "  -5.69%  "
$ws.Range("D35").Value = "3.599"
$ws.Range("E35").Value = "  -11.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.10"
$ws.Range("E36").Value = "  +6.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2140"
$ws.Range("E37").Value = "  -4.74%  "
$ws.Range("D38").Value = "5.084"
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("E39").Value = "  -5.51%  "
$ws.Range("D40").Value = "0.06045"
$ws.Range("E40").Value = "  -5.84%  "
$ws.Range("D41").Value = "0.6371"
$ws.Range("E41").Value = "  -5.68%  "
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("D43").Value = "8.008"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "1.407"
$ws.Range("E45").Value = "  -7.46%  "
$ws.Range("D46").Value = "13.66"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "3.718"
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("D48").Value = "0.5822"
$ws.Range("E48").Value = "  -6.18%  "
$ws.Range("D49").Value = "125.06"
$ws.Range("E49").Value = "  -3.93%  "
$ws.Range("E50").Value = "  -5.70%  "
$ws.Range("D51").Value = "0.06859"
$ws.Range("E51").Value = "  -4.12%  "
